$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to their rounded (2-decimal) equivalents.
$row5 = @{
    "B5" = 5.16
    "C5" = 3.81
    "D5" = 0.69
    "E5" = 11.47
    "F5" = 8.98
    "G5" = 3.83
    "H5" = 22.08
    "I5" = 6.33
    "J5" = 2.84
    "K5" = 4
    "L5" = 4.54
    "M5" = 4.99
    "N5" = 1.49
    "O5" = 4.14
    "P5" = 5.79
    "Q5" = 3.7
    "R5" = 0.37
    "S5" = 0.31
    "T5" = 55.62
    "U5" = 11.89
    "V5" = 3.82
    "W5" = 7.8
    "X5" = 4.05
    "Y5" = 0.53
    "Z5" = 10.76
    "AA5" = 3.37
    "AB5" = 3.1
    "AC5" = 3.6
    "AD5" = 4.88
    "AE5" = 0.48
    "AF5" = 20.72
    "AG5" = 2.03
    "AH5" = 4.77
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove row 6 entirely (shrinks the used range down to row 5).
$ws.Rows.Item(6).Delete()
